$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1) First table: add the 2014 row (row 24). Row 24 is currently completely
#    blank (the sheet jumps from row 23 straight to row 27), so no shifting
#    is required - we just need to populate it and copy the number formats
#    used by the rows above it.
# ---------------------------------------------------------------------------
$ws.Range("A23:J23").Copy()
$ws.Range("A24:J24").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A24").Value = 2014
$ws.Range("B24").Value = 163
$ws.Range("C24").Value = 3566
$ws.Range("D24").Value = 187
$ws.Range("E24").Value = 2648
$ws.Range("F24").Value = 464
$ws.Range("G24").Value = 220
$ws.Range("H24").Value = 2051
$ws.Range("I24").Value = 993
$ws.Range("J24").Value = 10292
$ws.Range("L24").Formula = "=B24+C24+F24+G24+H24"
$ws.Range("O24").Formula = "=J24+Z49"

# ---------------------------------------------------------------------------
# 2) Second table: add the 2014 row (row 49), also currently blank.
# ---------------------------------------------------------------------------
$ws.Range("A48:J48").Copy()
$ws.Range("A49:J49").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("Q48:Z48").Copy()
$ws.Range("Q49:Z49").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A49").Value = 2014
$ws.Range("B49").Value = 9
$ws.Range("C49").Value = 25
$ws.Range("D49").Value = 2
$ws.Range("E49").Formula = "=ROUND(E24*`$O49, 0)"
$ws.Range("F49").Value = 45
$ws.Range("G49").Value = 14
$ws.Range("H49").Value = 146
$ws.Range("I49").Value = 159
$ws.Range("J49").Value = 400
$ws.Range("L49").Formula = "=R49+S49+V49+W49+X49"
$ws.Range("O49").Formula = "=L49/L24"

$ws.Range("Q49").Value = 2014
$ws.Range("R49").Value = 9
$ws.Range("S49").Value = 25
$ws.Range("T49").Value = 2
$ws.Range("U49").Value = "NN"
$ws.Range("V49").Value = 45
$ws.Range("W49").Value = 14
$ws.Range("X49").Value = 146
$ws.Range("Y49").Value = 159
$ws.Range("Z49").Value = 400

# ---------------------------------------------------------------------------
# 3) Insert three rows before the old row 50 ("Total" label) so that it (and
#    everything under it) moves down to make room for a new row 51 plus two
#    blank spacer rows (matching the existing blank-row convention used
#    elsewhere in the sheet, e.g. rows 24-26, 49).
# ---------------------------------------------------------------------------
$ws.Rows.Item(50).Resize(3).Insert()

# New row 51: plain (unformatted) copy of the Q:Z block for 2014.
$ws.Range("Q51").Value = 2014
$ws.Range("R51").Value = 9
$ws.Range("S51").Value = 25
$ws.Range("T51").Value = 2
$ws.Range("U51").Value = "NN"
$ws.Range("V51").Value = 45
$ws.Range("W51").Value = 14
$ws.Range("X51").Value = 146
$ws.Range("Y51").Value = 159
$ws.Range("Z51").Value = 400

# ---------------------------------------------------------------------------
# 4) New row 74 (2014 total row at the bottom of the third table), matching
#    the formula pattern of the row above it (now row 73, previously 70).
# ---------------------------------------------------------------------------
$ws.Range("A73:J73").Copy()
$ws.Range("A74:J74").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A74").Value = 2014
$ws.Range("B74").Formula = "=B24+R49"
$ws.Range("C74").Formula = "=C24+S49"
$ws.Range("D74").Formula = "=D24+T49"
$ws.Range("E74").Formula = "=E24"
$ws.Range("F74").Formula = "=F24+V49"
$ws.Range("G74").Formula = "=G24+W49"
$ws.Range("H74").Formula = "=H24+X49"
$ws.Range("I74").Formula = "=I24+Y49"
$ws.Range("J74").Formula = "=J24+Z49"

# ---------------------------------------------------------------------------
# 5) Two pre-existing formulas that end up re-typed (value unchanged).
# ---------------------------------------------------------------------------
$ws.Range("O6").Formula = "=J6+Z31"
$ws.Range("L7").Formula = "=B7+C7+F7+G7+H7"

# ---------------------------------------------------------------------------
# 6) Sheet view bookkeeping, matching the diff.
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 40
$win.ScrollColumn = 1
[void]$ws.Range("A75").Select()
